$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (rows 2-4) duplicates column F's values
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("G4").Value = $ws.Range("F4").Value2
